# Update gh-pages to output generated at 456a3b4
# Applies updated "想去人数" (interest count) values to the F column
# across the 展览 / 演出 / 本地生活 / 全部类型 worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 1117
$ws1.Range("F5").Value  = 2756
$ws1.Range("F7").Value  = 689
$ws1.Range("F8").Value  = 64
$ws1.Range("F9").Value  = 266
$ws1.Range("F11").Value = 697
$ws1.Range("F12").Value = 101
$ws1.Range("F13").Value = 129
$ws1.Range("F14").Value = 1606
$ws1.Range("F17").Value = 199
$ws1.Range("F18").Value = 256

# --- Sheet: 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F10").Value = 23
$ws2.Range("F12").Value = 48

# --- Sheet: 本地生活 (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6352
$ws3.Range("F3").Value = 797
$ws3.Range("F5").Value = 255

# --- Sheet: 全部类型 (All Types, aggregated view) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 6352
$ws4.Range("F3").Value  = 797
$ws4.Range("F5").Value  = 255
$ws4.Range("F12").Value = 1117
$ws4.Range("F16").Value = 2756
$ws4.Range("F19").Value = 23
$ws4.Range("F21").Value = 48
$ws4.Range("F22").Value = 689
$ws4.Range("F23").Value = 64
$ws4.Range("F24").Value = 266
$ws4.Range("F27").Value = 697
$ws4.Range("F28").Value = 101
$ws4.Range("F29").Value = 129
$ws4.Range("F31").Value = 1606
$ws4.Range("F36").Value = 199
$ws4.Range("F43").Value = 256
